# feat: add 2022-Q4 data
#
# Plan:
#  - "总计" (summary) sheet: row 2 becomes the new 2022-Q4 totals, and the
#    old 2022-Q3 totals (previously row 2) are pushed down to a new row 3.
#  - The existing "2022-Q3" worksheet is duplicated (so its fund-holdings
#    data is preserved on the copy), the original slot is renamed to
#    "2022-Q4" and repopulated with the new quarter's fund-holdings data,
#    and the duplicate is renamed back to "2022-Q3" so it keeps the old data.

function Set-TextValue {
    # Forces a literal/text cell value (avoids the host auto-converting
    # numeric-looking strings like "6.42" or "233009" into numbers), while
    # leaving the cell on the default/unstyled format afterwards.
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "总计" sheet: shift the old 2022-Q3 row down to row 3, and write the
#    new 2022-Q4 figures into row 2.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)

# Give row 3 the same formatting as row 2 (style "2" on column A) before
# filling in the (old) 2022-Q3 values that used to live in row 2.
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q3"
$wsTotal.Range("C3").Value = 2
$wsTotal.Range("D3").Value = 0

# Row 2 now becomes the new 2022-Q4 summary line.
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.08

# ---------------------------------------------------------------------
# 2. Duplicate the "2022-Q3" sheet so the old fund-holdings data survives
#    on its own tab, then repurpose the original tab for 2022-Q4 data.
# ---------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Item(2)
$wsQ3.Copy($null, $wsQ3)

$wsQ3.Name = "2022-Q4"
$wsQ4Copy = $wb.Worksheets.Item(3)
$wsQ4Copy.Name = "2022-Q3"

$wsQ4 = $wb.Worksheets.Item(2)

# Reset page margins on the new 2022-Q4 tab to the workbook defaults
# (matching the "总计" sheet) instead of the ones inherited from the copy.
$wsQ4.PageSetup.LeftMargin = 54
$wsQ4.PageSetup.RightMargin = 54
$wsQ4.PageSetup.TopMargin = 72
$wsQ4.PageSetup.BottomMargin = 72
$wsQ4.PageSetup.HeaderMargin = 36
$wsQ4.PageSetup.FooterMargin = 36

# Clear the inherited 2022-Q3 fund-holdings data/formatting.
$wsQ4.Range("A1:H3").Clear()

# Re-create the header row with the "总计"-style formatting (cellXfs #2).
$wsTotal.Range("B1").Copy()
$wsQ4.Range("B1:H1").PasteSpecial(-4122)
$wsTotal.Range("A2").Copy()
$wsQ4.Range("A2:A3").PasteSpecial(-4122)

$wsQ4.Range("B1").Value = "基金代码"
$wsQ4.Range("C1").Value = "基金名称"
$wsQ4.Range("D1").Value = "基金规模"
$wsQ4.Range("E1").Value = "股票总仓位"
$wsQ4.Range("F1").Value = "仓位占比"
$wsQ4.Range("G1").Value = "持有市值(亿元)"
$wsQ4.Range("H1").Value = "仓位排名"

$wsQ4.Range("A2").Value = 0
Set-TextValue $wsQ4.Range("B2") "233009"
Set-TextValue $wsQ4.Range("C2") "大摩多因子精选策略混合"
Set-TextValue $wsQ4.Range("D2") "6.42"
Set-TextValue $wsQ4.Range("E2") "91.11"
Set-TextValue $wsQ4.Range("F2") "0.99"
Set-TextValue $wsQ4.Range("G2") "0.0636"
$wsQ4.Range("H2").Value = 3

$wsQ4.Range("A3").Value = 1
Set-TextValue $wsQ4.Range("B3") "512190"
Set-TextValue $wsQ4.Range("C3") "浙商汇金中证浙江凤凰行动50ETF"
Set-TextValue $wsQ4.Range("D3") "0.50"
Set-TextValue $wsQ4.Range("E3") "99.24"
Set-TextValue $wsQ4.Range("F3") "3.48"
Set-TextValue $wsQ4.Range("G3") "0.0174"
$wsQ4.Range("H3").Value = 8
